# ------------------------------------------------------------------
# Re-layout the "history" worksheet to match the updated data schema:
#   - the "antenne" column moves from position I to the end (col X)
#   - the "products" column moves from position V to the very end (col Y)
#   - every field that used to sit between "antenne" and "products"
#     shifts one column to the left to fill the gap
#   - a new history record (row 5) is appended
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-NumberCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# --- Update header row (row 1) to the new column layout ---
Set-TextCell "A1" "id"
Set-TextCell "B1" "timestamp"
Set-TextCell "C1" "documentNumber"
Set-TextCell "D1" "date"
Set-TextCell "E1" "dateDepart"
Set-TextCell "F1" "clientId"
Set-TextCell "G1" "clientName"
Set-TextCell "H1" "destination"
Set-TextCell "I1" "itineraire"
Set-TextCell "J1" "driverId"
Set-TextCell "K1" "driverName"
Set-TextCell "L1" "driverCIN"
Set-TextCell "M1" "driverPhone"
Set-TextCell "N1" "vehicleMatricule"
Set-TextCell "O1" "vehicleModel"
Set-TextCell "P1" "convoyeurId"
Set-TextCell "Q1" "convoyeurName"
Set-TextCell "R1" "convoyeurCard"
Set-TextCell "S1" "convoyeurCIN"
Set-TextCell "T1" "convoyeurPhone"
Set-TextCell "U1" "passavantNumber"
Set-TextCell "V1" "passavantExpiry"
Set-TextCell "W1" "bonLivraison"
Set-TextCell "X1" "antenne"
Set-TextCell "Y1" "products"

# --- Rewrite existing data rows (2-4) into the new column layout ---
# Row 2
Set-NumberCell "A2" 1766599794753
Set-TextCell "B2" "2025-12-24T18:09:54.753Z"
Set-TextCell "C2" "10/2025"
Set-TextCell "D2" "2025-12-24"
Set-TextCell "E2" "2025-12-24T07:09"
Set-NumberCell "F2" 2
Set-TextCell "G2" "Client B"
Set-TextCell "H2" "Warehouse B"
Set-TextCell "I2" "Point D, Point E, Point F"
Set-NumberCell "J2" 4
Set-TextCell "K2" "oudra"
Set-TextCell "L2" "p3608237"
Set-TextCell "M2" "0567233893"
Set-TextCell "N2" "sdgas556"
Set-TextCell "O2" "volvo"
Set-NumberCell "P2" 2
Set-TextCell "Q2" "Youssef Alaoui"
Set-TextCell "R2" "CCE002"
Set-TextCell "S2" "YA333444"
Set-TextCell "T2" "0656789012"
Set-TextCell "U2" "23423432"
Set-TextCell "V2" "2025-12-25"
Set-TextCell "W2" "43324"
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()

# Row 3
Set-NumberCell "A3" 1766597982396
Set-TextCell "B3" "2025-12-24T17:39:42.396Z"
Set-TextCell "C3" "6/2025"
Set-TextCell "D3" "2025-12-25"
Set-TextCell "E3" "2025-12-06T18:21"
Set-NumberCell "F3" 4
Set-TextCell "G3" "ABC Company"
Set-TextCell "H3" "ABC Warehouse"
Set-TextCell "I3" "Route 1, Route 2"
Set-NumberCell "J3" 4
Set-TextCell "K3" "oudra"
Set-TextCell "L3" "p3608237"
Set-TextCell "M3" "0567233893"
Set-TextCell "N3" "sdgas556"
Set-TextCell "O3" "volvo"
Set-NumberCell "P3" 1
Set-TextCell "Q3" "Omar Tazi"
Set-TextCell "R3" "CCE001"
Set-TextCell "S3" "OT111222"
Set-TextCell "T3" "0645678901"
Set-TextCell "U3" "234234"
Set-TextCell "V3" "2025-12-26"
Set-TextCell "W3" "34324"
Set-TextCell "X3" "Antenne 3"
$ws.Range("Y3").ClearContents()

# Row 4
Set-NumberCell "A4" 1766596896691
Set-TextCell "B4" "2025-12-24T17:21:36.691Z"
Set-TextCell "C4" "6/2025"
Set-TextCell "D4" "2025-12-25"
Set-TextCell "E4" "2025-12-06T18:21"
Set-NumberCell "F4" 4
Set-TextCell "G4" "ABC Company"
Set-TextCell "H4" "ABC Warehouse"
Set-TextCell "I4" "Route 1, Route 2"
Set-NumberCell "J4" 1
Set-TextCell "K4" "Ahmed Benali"
Set-TextCell "L4" "AB123456"
Set-TextCell "M4" "0612345678"
Set-TextCell "N4" ""
Set-TextCell "O4" ""
Set-NumberCell "P4" 1
Set-TextCell "Q4" "Omar Tazi"
Set-TextCell "R4" "CCE001"
Set-TextCell "S4" "OT111222"
Set-TextCell "T4" "0645678901"
Set-TextCell "U4" "234234"
Set-TextCell "V4" "2025-12-26"
Set-TextCell "W4" "34324"
Set-TextCell "X4" "Antenne 3"
$ws.Range("Y4").ClearContents()

# --- Append new history row (row 5) ---
Set-NumberCell "A5" 1766961828533
Set-TextCell "B5" "2025-12-28T22:43:48.533Z"
Set-TextCell "C5" "26/2025"
Set-TextCell "D5" "2025-12-03"
Set-TextCell "E5" "2025-12-19T23:43"
Set-NumberCell "F5" 1
Set-TextCell "G5" "SFI"
Set-TextCell "H5" "SFI Depot"
Set-TextCell "I5" "Point A, Point B, Point C"
Set-NumberCell "J5" 1
Set-TextCell "K5" "Ahmed Benali"
Set-TextCell "L5" "AB123456"
Set-TextCell "M5" "0612345678"
Set-TextCell "N5" ""
Set-TextCell "O5" ""
Set-NumberCell "P5" 1
Set-TextCell "Q5" "Omar Tazi"
Set-TextCell "R5" "CCE001"
Set-TextCell "S5" "OT111222"
Set-TextCell "T5" "0645678901"
Set-TextCell "U5" "1221"
Set-TextCell "V5" "2025-12-12"
Set-TextCell "W5" "121212"
$ws.Range("X5").ClearContents()  # no antenne value for this row
Set-TextCell "Y5" "[{`"name`":`"Produit A`",`"quantity`":`"410`",`"unit`":`"Kg`"}]"
